# Adapt column header formatting to respective input file names (#7)
# - Rename header suffixes "_old" -> "_FV2310" and "_new" -> "_FV2404"
# - Turn A1:U66 into an Excel Table ("Table1") with the renamed headers
# - Freeze the header row (top row) in the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row (row 1) cell values.
# ---------------------------------------------------------------------------
$headers = @(
  "Segmentname_FV2310",
  "Segmentgruppe_FV2310",
  "Segment_FV2310",
  "Datenelement_FV2310",
  "Segment ID_FV2310",
  "Code_FV2310",
  "Qualifier_FV2310",
  "Beschreibung_FV2310",
  "Bedingungsausdruck_FV2310",
  "Bedingung_FV2310",
  "diff",
  "Segmentname_FV2404",
  "Segmentgruppe_FV2404",
  "Segment_FV2404",
  "Datenelement_FV2404",
  "Segment ID_FV2404",
  "Code_FV2404",
  "Qualifier_FV2404",
  "Beschreibung_FV2404",
  "Bedingungsausdruck_FV2404",
  "Bedingung_FV2404"
)

$headerRange = $ws.Range("A1:U1")

for ($i = 0; $i -lt $headers.Count; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------------
# 2) Turn A1:U66 into a native Excel Table ("ListObject").
#    Clear the pre-existing header formatting first so the table creation
#    doesn't capture it into a one-off header dxf/style; then restore the
#    original header look (bold, centered, wrapped, filled, bordered)
#    directly on the cells afterwards.
# ---------------------------------------------------------------------------
$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U66")
$tbl = $ws.ListObjects.Add(
  [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
  $tableRange,
  $null,
  [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
$headerRange.WrapText = $true
$headerRange.Interior.Color = 14277081
$headerRange.Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$headerRange.Borders.Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin

# ---------------------------------------------------------------------------
# 3) Freeze the header row.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Header renamed, table '$($tbl.Name)' created, top row frozen."
